$wb = $excel.ActiveWorkbook

# Helper: write a value to a cell as TEXT. Excel's COM `.Value` setter
# auto-coerces strings that look like numbers or booleans ("3", "True",
# "False", ...) into Number/Boolean cells. The source data must stay text
# (matching t="inlineStr" in the target file), so any value that would be
# misread is written with a leading apostrophe, which forces Excel to keep
# it as a literal text value. Plain, unambiguous strings are set directly
# so no extra formatting is introduced.
function Set-TextValue($cell, $val) {
    $needsEscape = $false
    if ($val -eq "True" -or $val -eq "False") {
        $needsEscape = $true
    } elseif ($val -match '^[0-9]+$') {
        $needsEscape = $true
    }
    if ($needsEscape) {
        $cell.Value = "'" + $val
    } else {
        $cell.Value = $val
    }
}

# ---------------------------------------------------------------------------
# Sheet "Ares Condicionados": replace the single existing data row (row 2)
# with six new data rows (rows 2-7). Columns: A=Comodo, B=Nome, C=Ligado
# (all text), D=Temperatura, E=Intensidade (numeric).
# ---------------------------------------------------------------------------
$wsAr = $wb.Worksheets.Item("Ares Condicionados")

$arRows = @(
    @{ A = "Bom";   B = "Ar1";  C = "False"; D = 0;  E = 0  },
    @{ A = "Bom";   B = "3";    C = "True";  D = 15; E = 53 },
    @{ A = "Bom";   B = "4";    C = "False"; D = 0;  E = 0  },
    @{ A = "sorte"; B = "1";    C = "True";  D = 13; E = 57 },
    @{ A = "sorte"; B = "23";   C = "False"; D = 0;  E = 0  },
    @{ A = "sorte"; B = "4312"; C = "False"; D = 0;  E = 0  }
)

$r = 2
foreach ($row in $arRows) {
    Set-TextValue $wsAr.Cells.Item($r, 1) $row.A
    Set-TextValue $wsAr.Cells.Item($r, 2) $row.B
    Set-TextValue $wsAr.Cells.Item($r, 3) $row.C
    $wsAr.Cells.Item($r, 4).Value = $row.D
    $wsAr.Cells.Item($r, 5).Value = $row.E
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "Janelas": add three new data rows (rows 2-4). Columns: A=Comodo,
# B=Nome (text), F=Abertura (numeric), G=Tranca (text True/False).
# ---------------------------------------------------------------------------
$wsJa = $wb.Worksheets.Item("Janelas")

$jaRows = @(
    @{ A = "sorte"; B = "fdsoj"; F = 0; G = "True"  },
    @{ A = "Bom";   B = "1";     F = 0; G = "False" },
    @{ A = "sorte"; B = "2";     F = 0; G = "False" }
)

$r = 2
foreach ($row in $jaRows) {
    Set-TextValue $wsJa.Cells.Item($r, 1) $row.A
    Set-TextValue $wsJa.Cells.Item($r, 2) $row.B
    $wsJa.Cells.Item($r, 6).Value = $row.F
    Set-TextValue $wsJa.Cells.Item($r, 7) $row.G
    $r = $r + 1
}
